$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountSet")

# Row 3: Credit: Curr Stmt Bal
$ws.Cells.Item(3, 1).Value = "Credit: Curr Stmt Bal"
$ws.Cells.Item(3, 2).Value = 100
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 20000
$ws.Cells.Item(3, 5).Value = "curr stmt bal"

# Row 4: Credit: Prev Stmt Bal
$ws.Cells.Item(4, 1).Value = "Credit: Prev Stmt Bal"
$ws.Cells.Item(4, 2).Value = 100
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 20000
$ws.Cells.Item(4, 5).Value = "prev stmt bal"
$ws.Cells.Item(4, 6).Value = "'20000102"
$ws.Cells.Item(4, 7).Value = "compound"
$ws.Cells.Item(4, 8).Value = 0.01
$ws.Cells.Item(4, 9).Value = "monthly"
$ws.Cells.Item(4, 10).Value = 40

# Row 5: test loan: Principal Balance
$ws.Cells.Item(5, 1).Value = "test loan: Principal Balance"
$ws.Cells.Item(5, 2).Value = 100
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 9999
$ws.Cells.Item(5, 5).Value = "principal balance"
$ws.Cells.Item(5, 6).Value = "'20000102"
$ws.Cells.Item(5, 7).Value = "simple"
$ws.Cells.Item(5, 8).Value = 0.01
$ws.Cells.Item(5, 9).Value = "daily"
$ws.Cells.Item(5, 10).Value = 50

# Row 6: test loan: Interest
$ws.Cells.Item(6, 1).Value = "test loan: Interest"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 9999
$ws.Cells.Item(6, 5).Value = "interest"
